$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: ECs -> Gnai2 -> Adcy1 -> FAPs
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Gnai2"
$ws.Range("C2").Value = "Adcy1"
$ws.Range("D2").Value = "FAPs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 203.7816646666667
$ws.Range("H2").Value = 611.344994
$ws.Range("I2").Value = 0.6667327591988204
$ws.Range("J2").Value = 0.6667327591988205
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.06943366666666666
$ws.Range("N2").Value = 0.208301
$ws.Range("O2").Value = 0.2790002116266049
$ws.Range("P2").Value = 0.2790002116266049
$ws.Range("Q2").Value = 14.14930817724378
$ws.Range("R2").Value = 127.343773595194
$ws.Range("S2").Value = 0.1860185809148611
$ws.Range("T2").Value = 0.1860185809148611

# Row 3: ECs -> Gnai2 -> Adcy1 -> sCs
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Gnai2"
$ws.Range("C3").Value = "Adcy1"
$ws.Range("D3").Value = "sCs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 203.7816646666667
$ws.Range("H3").Value = 611.344994
$ws.Range("I3").Value = 0.6667327591988204
$ws.Range("J3").Value = 0.6667327591988205
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.1794323333333333
$ws.Range("N3").Value = 0.538297
$ws.Range("O3").Value = 0.7209997883733951
$ws.Range("P3").Value = 0.7209997883733951
$ws.Range("Q3").Value = 36.56501958169089
$ws.Range("R3").Value = 329.0851762352181
$ws.Range("S3").Value = 0.4807141782839592
$ws.Range("T3").Value = 0.4807141782839593

# Row 4: FAPs -> Gnai2 -> Adcy1 -> FAPs
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Gnai2"
$ws.Range("C4").Value = "Adcy1"
$ws.Range("D4").Value = "FAPs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 63.14058933333333
$ws.Range("H4").Value = 189.421768
$ws.Range("I4").Value = 0.2065833519051582
$ws.Range("J4").Value = 0.2065833519051582
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.06943366666666666
$ws.Range("N4").Value = 0.208301
$ws.Range("O4").Value = 0.2790002116266049
$ws.Range("P4").Value = 0.2790002116266049
$ws.Range("Q4").Value = 4.384082632907555
$ws.Range("R4").Value = 39.456743696168
$ws.Range("S4").Value = 0.05763679890007253
$ws.Range("T4").Value = 0.05763679890007253

# Row 5: FAPs -> Gnai2 -> Adcy1 -> sCs
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Gnai2"
$ws.Range("C5").Value = "Adcy1"
$ws.Range("D5").Value = "sCs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 63.14058933333333
$ws.Range("H5").Value = 189.421768
$ws.Range("I5").Value = 0.2065833519051582
$ws.Range("J5").Value = 0.2065833519051582
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.1794323333333333
$ws.Range("N5").Value = 0.538297
$ws.Range("O5").Value = 0.7209997883733951
$ws.Range("P5").Value = 0.7209997883733951
$ws.Range("Q5").Value = 11.32946327212178
$ws.Range("R5").Value = 101.965169449096
$ws.Range("S5").Value = 0.1489465530050856
$ws.Range("T5").Value = 0.1489465530050857

# Row 6: sCs -> Gnai2 -> Adcy1 -> FAPs
$ws.Range("A6").Value = "sCs"
$ws.Range("B6").Value = "Gnai2"
$ws.Range("C6").Value = "Adcy1"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 38.719942
$ws.Range("H6").Value = 116.159826
$ws.Range("I6").Value = 0.1266838888960214
$ws.Range("J6").Value = 0.1266838888960214
$ws.Range("K6").Value = 1
$ws.Range("L6").Value = 0.3333333333333333
$ws.Range("M6").Value = 0.06943366666666666
$ws.Range("N6").Value = 0.208301
$ws.Range("O6").Value = 0.2790002116266049
$ws.Range("P6").Value = 0.2790002116266049
$ws.Range("Q6").Value = 2.688467546180666
$ws.Range("R6").Value = 24.196207915626
$ws.Range("S6").Value = 0.03534483181167127
$ws.Range("T6").Value = 0.03534483181167128

# Row 7: sCs -> Gnai2 -> Adcy1 -> sCs
$ws.Range("A7").Value = "sCs"
$ws.Range("B7").Value = "Gnai2"
$ws.Range("C7").Value = "Adcy1"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 38.719942
$ws.Range("H7").Value = 116.159826
$ws.Range("I7").Value = 0.1266838888960214
$ws.Range("J7").Value = 0.1266838888960214
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.1794323333333333
$ws.Range("N7").Value = 0.538297
$ws.Range("O7").Value = 0.7209997883733951
$ws.Range("P7").Value = 0.7209997883733951
$ws.Range("Q7").Value = 6.947609539591332
$ws.Range("R7").Value = 62.528485856322
$ws.Range("S7").Value = 0.09133905708435011
$ws.Range("T7").Value = 0.09133905708435014
